$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.113.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.975.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4971"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4221"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09307"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.966.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.925"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.475"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001111"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06706"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.980"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.158.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.205.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.349"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.272"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09867"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.518"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.737"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02432"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.331"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06457"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.061"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6502"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2011"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.009"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6231"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.366"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.192"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.483"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000327"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06971"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
